$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.736.80"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.914.83"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "2.913.38"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "3.396.50"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "62.643.53"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "2.911.72"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.662"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("E29").Value = "  +5.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "2.720.74"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0340"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "354.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000223"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +17.54%  "
$ws.Range("E51").Value = "  -0.22%  "
